# Auto-generated script applying the Maduin_Profits.xlsx data refresh diff
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("M21").ClearContents()
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("M23").ClearContents()
$ws.Range("H74").Value = 3500
$ws.Range("H77").Value = 3500
$ws.Range("H137").Value = 1591.5
$ws.Range("I137").Value = 1521.6154
$ws.Range("J137").Value = 2500
$ws.Range("K137").Value = 4564.8462
$ws.Range("L137").Value = 7500
$ws.Range("M137").Value = -2014.8462
$ws.Range("H138").Value = 3838.3333
$ws.Range("I138").Value = 8000
$ws.Range("J138").Value = 3630.25
$ws.Range("K138").Value = 24000
$ws.Range("L138").Value = 10890.75
$ws.Range("M138").Value = -18860
$ws.Range("N138").Value = -21170.75

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1100
$ws.Range("I2").Value = 1100
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1100
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -987
$ws.Range("H92").Value = 42412.25
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 42412.25
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 42412.25
$ws.Range("N92").Value = -47404.25
$ws.Range("H94").Value = 38999.5
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 38999.5
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 38999.5
$ws.Range("N94").Value = -40801.5
$ws.Range("H97").Value = 12160.6
$ws.Range("I97").Value = 262.5
$ws.Range("J97").Value = 20092.666
$ws.Range("K97").Value = 262.5
$ws.Range("L97").Value = 20092.666
$ws.Range("M97").Value = 233.5
$ws.Range("N97").Value = -21084.666
$ws.Range("H116").Value = 1100
$ws.Range("I116").Value = 1100
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 1100
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 1194

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1100
$ws.Range("I3").Value = 1100
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1100
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -986
$ws.Range("H37").Value = 1025.8
$ws.Range("I37").Value = 1250
$ws.Range("J37").Value = 129
$ws.Range("K37").Value = 1250
$ws.Range("L37").Value = 129
$ws.Range("M37").Value = -1113
$ws.Range("H94").Value = 1968.1666
$ws.Range("I94").Value = 1809
$ws.Range("J94").Value = 2000
$ws.Range("K94").Value = 1809
$ws.Range("L94").Value = 2000
$ws.Range("M94").Value = -1358
$ws.Range("N94").Value = -2902
$ws.Range("H99").Value = 2387.5
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 2387.5
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 2387.5
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -5383.5
$ws.Range("H107").Value = 629.625
$ws.Range("I107").Value = 660.5714
$ws.Range("J107").Value = 413
$ws.Range("K107").Value = 660.5714
$ws.Range("L107").Value = 413
$ws.Range("M107").Value = 1259.4286
$ws.Range("H134").Value = 4201.8335
$ws.Range("I134").Value = 1105.5
$ws.Range("J134").Value = 5750
$ws.Range("K134").Value = 3316.5
$ws.Range("L134").Value = 17250
$ws.Range("M134").Value = -781.5
$ws.Range("N134").Value = -22320

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 77
$ws.Range("I2").Value = 77
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 77
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 36
$ws.Range("N2").ClearContents()
$ws.Range("H107").Value = 584.25
$ws.Range("I107").Value = 377.375
$ws.Range("J107").Value = 998
$ws.Range("K107").Value = 377.375
$ws.Range("L107").Value = 998
$ws.Range("M107").Value = 1542.625
$ws.Range("N107").Value = -4838
$ws.Range("H119").Value = 60000
$ws.Range("I119").Value = 0
$ws.Range("J119").Value = 60000
$ws.Range("K119").Value = 0
$ws.Range("L119").Value = 60000
$ws.Range("N119").Value = -69676

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H64").Value = 142.5
$ws.Range("I64").Value = 142.5
$ws.Range("J64").Value = 0
$ws.Range("K64").Value = 427.5
$ws.Range("L64").Value = 0
$ws.Range("M64").Value = -157.5
$ws.Range("H67").Value = 142.5
$ws.Range("I67").Value = 142.5
$ws.Range("J67").Value = 0
$ws.Range("K67").Value = 427.5
$ws.Range("L67").Value = 0
$ws.Range("M67").Value = 508.5
$ws.Range("H140").Value = 1745.875
$ws.Range("I140").Value = 665.5
$ws.Range("J140").Value = 4987
$ws.Range("K140").Value = 1996.5
$ws.Range("L140").Value = 14961
$ws.Range("M140").Value = 3183.5
$ws.Range("N140").Value = -25321

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 733.3125
$ws.Range("I2").Value = 1018.1
$ws.Range("J2").Value = 258.66666
$ws.Range("K2").Value = 1018.1
$ws.Range("L2").Value = 258.66666
$ws.Range("M2").Value = -905.1
$ws.Range("H46").Value = 11660.125
$ws.Range("I46").Value = 4466.6665
$ws.Range("J46").Value = 15976.2
$ws.Range("K46").Value = 4466.6665
$ws.Range("L46").Value = 15976.2
$ws.Range("M46").Value = -4310.6665
$ws.Range("H52").Value = 35596.92
$ws.Range("H70").Value = 7611.75
$ws.Range("I70").Value = 9974.25
$ws.Range("J70").Value = 5249.25
$ws.Range("K70").Value = 9974.25
$ws.Range("L70").Value = 5249.25
$ws.Range("M70").Value = -9704.25
$ws.Range("N70").Value = -5789.25
$ws.Range("H73").Value = 7611.75
$ws.Range("I73").Value = 9974.25
$ws.Range("J73").Value = 5249.25
$ws.Range("K73").Value = 9974.25
$ws.Range("L73").Value = 5249.25
$ws.Range("M73").Value = -9038.25
$ws.Range("N73").Value = -7121.25
$ws.Range("H92").Value = 12876.857
$ws.Range("I92").Value = 0
$ws.Range("J92").Value = 12876.857
$ws.Range("K92").Value = 0
$ws.Range("L92").Value = 12876.857
$ws.Range("N92").Value = -16620.857
$ws.Range("H122").Value = 1477.3334
$ws.Range("I122").Value = 1477.3334
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 4432.0002
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -1982.0002

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1700.0625
$ws.Range("I46").Value = 1540.2
$ws.Range("J46").Value = 1772.7273
$ws.Range("K46").Value = 1540.2
$ws.Range("L46").Value = 1772.7273
$ws.Range("M46").Value = -1352.2
$ws.Range("H68").Value = 1661.3334
$ws.Range("I68").Value = 997
$ws.Range("J68").Value = 2990
$ws.Range("K68").Value = 997
$ws.Range("L68").Value = 2990
$ws.Range("M68").Value = -248
$ws.Range("N68").Value = -4488
$ws.Range("H71").Value = 1661.3334
$ws.Range("I71").Value = 997
$ws.Range("J71").Value = 2990
$ws.Range("K71").Value = 4985
$ws.Range("L71").Value = 14950
$ws.Range("M71").Value = -1241
$ws.Range("N71").Value = -22438
$ws.Range("H124").Value = 0
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 24500000
$ws.Range("I3").Value = 4000000
$ws.Range("J3").Value = 45000000
$ws.Range("K3").Value = 4000000
$ws.Range("L3").Value = 45000000
$ws.Range("M3").Value = -3999886
$ws.Range("N3").Value = -45000228
$ws.Range("H122").Value = 844.6667
$ws.Range("I122").Value = 764.5
$ws.Range("J122").Value = 1005
$ws.Range("K122").Value = 2293.5
$ws.Range("L122").Value = 3015
$ws.Range("M122").Value = 156.5
$ws.Range("N122").Value = -7915
$ws.Range("H132").Value = 167544.83
$ws.Range("I132").Value = 333589.66
$ws.Range("J132").Value = 1500
$ws.Range("K132").Value = 1000768.98
$ws.Range("L132").Value = 4500
$ws.Range("M132").Value = -998238.98
$ws.Range("N132").Value = -9560
